# Update "想去人数" (want-to-go count) figures in column F across the
# three sheets that carry this data: 展览 (sheet 1), 演出 (sheet 2) and
# 全部类型 (sheet 4, the combined view). 本地生活 (sheet 3) is unaffected.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# 展览
$ws1.Range("F2").Value = 7617
$ws1.Range("F3").Value = 7617
$ws1.Range("F5").Value = 7807
$ws1.Range("F9").Value = 6517
$ws1.Range("F10").Value = 3341
$ws1.Range("F12").Value = 3698
$ws1.Range("F14").Value = 39
$ws1.Range("F16").Value = 58
$ws1.Range("F17").Value = 45
$ws1.Range("F18").Value = 460
$ws1.Range("F20").Value = 307
$ws1.Range("F21").Value = 320
$ws1.Range("F22").Value = 3792
$ws1.Range("F25").Value = 953
$ws1.Range("F27").Value = 1440
$ws1.Range("F28").Value = 76
$ws1.Range("F30").Value = 2723
$ws1.Range("F31").Value = 1750
$ws1.Range("F35").Value = 3579
$ws1.Range("F36").Value = 281
$ws1.Range("F41").Value = 1383
$ws1.Range("F42").Value = 241
$ws1.Range("F43").Value = 542
$ws1.Range("F44").Value = 629

# 演出
$ws2.Range("F9").Value = 101
$ws2.Range("F13").Value = 86
$ws2.Range("F16").Value = 1

# 全部类型
$ws4.Range("F5").Value = 7617
$ws4.Range("F6").Value = 7617
$ws4.Range("F8").Value = 7807
$ws4.Range("F11").Value = 6517
$ws4.Range("F12").Value = 3341
$ws4.Range("F14").Value = 3698
$ws4.Range("F16").Value = 39
$ws4.Range("F18").Value = 58
$ws4.Range("F19").Value = 45
$ws4.Range("F20").Value = 460
$ws4.Range("F21").Value = 307
$ws4.Range("F22").Value = 101
$ws4.Range("F23").Value = 320
$ws4.Range("F24").Value = 3792
$ws4.Range("F30").Value = 953
$ws4.Range("F32").Value = 1440
$ws4.Range("F33").Value = 76
$ws4.Range("F35").Value = 2723
$ws4.Range("F36").Value = 1750
$ws4.Range("F39").Value = 86
$ws4.Range("F40").Value = 3579
$ws4.Range("F41").Value = 281
$ws4.Range("F46").Value = 1383
$ws4.Range("F47").Value = 241
$ws4.Range("F49").Value = 542
$ws4.Range("F50").Value = 629
